# Rename the sheet from "List 1" to "DATA"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DATA"

# The header cells in C2:D2, C3, C4, C5:D5 carried a distinct
# (Calibri / light-gray-fill) character format left over from a
# paste. Normalize them back to the sheet's regular formatting by
# copying the format already used by the rest of the data rows
# (e.g. A2, which uses the plain/default look) onto them.
$normalFormat = $ws.Range("A2")

$normalFormat.Copy()
$ws.Range("C2:D2").PasteSpecial(-4122)

$normalFormat.Copy()
$ws.Range("C3").PasteSpecial(-4122)

$normalFormat.Copy()
$ws.Range("C4").PasteSpecial(-4122)

$normalFormat.Copy()
$ws.Range("C5:D5").PasteSpecial(-4122)

$excel.CutCopyMode = 0
